$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config")
$wf = $excel.WorksheetFunction
$members = $wf | Get-Member
$names = $members | ForEach-Object { $_.Name }
Write-Output ($names -join "|")
